$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix B20 to be a numeric value instead of a text/inline-string value
$ws.Cells.Item(20, 2).Value = 3

# Add new row 21 with data
$ws.Cells.Item(21, 1).Value = "Ying Tang"

# B21 must stay as text "1" (not be auto-converted to a number)
$ws.Cells.Item(21, 2).NumberFormat = "@"
$ws.Cells.Item(21, 2).Value = "1"
$ws.Cells.Item(21, 2).NumberFormat = "General"
$ws.Cells.Item(21, 2).ClearFormats()

$ws.Cells.Item(21, 3).Value = "No clear novelty"
$ws.Cells.Item(21, 4).Value = "CRT"
$ws.Cells.Item(21, 5).Value = "MET"
$ws.Cells.Item(21, 6).Value = "4efacd8b-a5d8-471d-9660-f5eb687b96fc"
$ws.Cells.Item(21, 7).Value = "Byni8NLHf_annotated.xlsx"
$ws.Cells.Item(21, 8).Value = "No clear novelty"
